$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 26.83824466666667
$ws.Range("H2").Value = 80.514734
$ws.Range("I2").Value = 0.8882651037973995
$ws.Range("J2").Value = 0.8882651037973996
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 73.19890333333333
$ws.Range("N2").Value = 219.59671
$ws.Range("O2").Value = 0.6596328743217019
$ws.Range("P2").Value = 0.6596328743217019
$ws.Range("Q2").Value = 1964.530076991682
$ws.Range("R2").Value = 17680.77069292514
$ws.Range("S2").Value = 0.5859288635775436
$ws.Range("T2").Value = 0.5859288635775436

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 26.83824466666667
$ws.Range("H3").Value = 80.514734
$ws.Range("I3").Value = 0.8882651037973995
$ws.Range("J3").Value = 0.8882651037973996
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.661646333333335
$ws.Range("N3").Value = 28.984939
$ws.Range("O3").Value = 0.08706605224007774
$ws.Range("P3").Value = 0.08706605224007773
$ws.Range("Q3").Value = 259.301628176803
$ws.Range("R3").Value = 2333.714653591227
$ws.Range("S3").Value = 0.07733773593026247
$ws.Range("T3").Value = 0.07733773593026246

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 26.83824466666667
$ws.Range("H4").Value = 80.514734
$ws.Range("I4").Value = 0.8882651037973995
$ws.Range("J4").Value = 0.8882651037973996
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.10860633333333
$ws.Range("N4").Value = 84.325819
$ws.Range("O4").Value = 0.2533010734382203
$ws.Range("P4").Value = 0.2533010734382204
$ws.Range("Q4").Value = 754.3856540130162
$ws.Range("R4").Value = 6789.470886117146
$ws.Range("S4").Value = 0.2249985042895935
$ws.Range("T4").Value = 0.2249985042895936

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.239011
$ws.Range("H5").Value = 3.717033
$ws.Range("I5").Value = 0.04100753414354395
$ws.Range("J5").Value = 0.04100753414354396
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 73.19890333333333
$ws.Range("N5").Value = 219.59671
$ws.Range("O5").Value = 0.6596328743217019
$ws.Range("P5").Value = 0.6596328743217019
$ws.Range("Q5").Value = 90.69424641793665
$ws.Range("R5").Value = 816.2482177614299
$ws.Range("S5").Value = 0.02704991761595123
$ws.Range("T5").Value = 0.02704991761595123

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.239011
$ws.Range("H6").Value = 3.717033
$ws.Range("I6").Value = 0.04100753414354395
$ws.Range("J6").Value = 0.04100753414354396
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.661646333333335
$ws.Range("N6").Value = 28.984939
$ws.Range("O6").Value = 0.08706605224007774
$ws.Range("P6").Value = 0.08706605224007773
$ws.Range("Q6").Value = 11.97088608510967
$ws.Range("R6").Value = 107.737974765987
$ws.Range("S6").Value = 0.00357036410997857
$ws.Range("T6").Value = 0.00357036410997857

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.239011
$ws.Range("H7").Value = 3.717033
$ws.Range("I7").Value = 0.04100753414354395
$ws.Range("J7").Value = 0.04100753414354396
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.10860633333333
$ws.Range("N7").Value = 84.325819
$ws.Range("O7").Value = 0.2533010734382203
$ws.Range("P7").Value = 0.2533010734382204
$ws.Range("Q7").Value = 34.82687244166966
$ws.Range("R7").Value = 313.441851975027
$ws.Range("S7").Value = 0.01038725241761415
$ws.Range("T7").Value = 0.01038725241761416

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.136972666666667
$ws.Range("H8").Value = 6.410918000000001
$ws.Range("I8").Value = 0.07072736205905639
$ws.Range("J8").Value = 0.0707273620590564
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 73.19890333333333
$ws.Range("N8").Value = 219.59671
$ws.Range("O8").Value = 0.6596328743217019
$ws.Range("P8").Value = 0.6596328743217019
$ws.Range("Q8").Value = 156.4240556533089
$ws.Range("R8").Value = 1407.81650087978
$ws.Range("S8").Value = 0.04665409312820705
$ws.Range("T8").Value = 0.04665409312820706

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.136972666666667
$ws.Range("H9").Value = 6.410918000000001
$ws.Range("I9").Value = 0.07072736205905639
$ws.Range("J9").Value = 0.0707273620590564
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.661646333333335
$ws.Range("N9").Value = 28.984939
$ws.Range("O9").Value = 0.08706605224007774
$ws.Range("P9").Value = 0.08706605224007773
$ws.Range("Q9").Value = 20.64667412933356
$ws.Range("R9").Value = 185.820067164002
$ws.Range("S9").Value = 0.006157952199836695
$ws.Range("T9").Value = 0.006157952199836696

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.136972666666667
$ws.Range("H10").Value = 6.410918000000001
$ws.Range("I10").Value = 0.07072736205905639
$ws.Range("J10").Value = 0.0707273620590564
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 28.10860633333333
$ws.Range("N10").Value = 84.325819
$ws.Range("O10").Value = 0.2533010734382203
$ws.Range("P10").Value = 0.2533010734382204
$ws.Range("Q10").Value = 60.06732343242689
$ws.Range("R10").Value = 540.605910891842
$ws.Range("S10").Value = 0.01791531673101264
$ws.Range("T10").Value = 0.01791531673101265

